# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.846.00"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.59"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7361"
$ws.Range("E5").Value = "  -5.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.79"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3155"
$ws.Range("E8").Value = "  +0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.76"
$ws.Range("E9").Value = "  -4.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07095"
$ws.Range("E10").Value = "  -2.44%  "

$ws.Range("E11").Value = "  -8.73%  "

$ws.Range("E12").Value = "  -3.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.417"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.866.55"
$ws.Range("E14").Value = "  -1.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.57"
$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.849.96"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.049"
$ws.Range("E17").Value = "  -3.08%  "

$ws.Range("E18").Value = "  -3.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.89"
$ws.Range("E19").Value = "  -1.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007833"
$ws.Range("E20").Value = "  -0.94%  "

$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.116.78"
$ws.Range("E22").Value = "  -2.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.902"
$ws.Range("E23").Value = "  -3.52%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1568"
$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.307"
$ws.Range("E26").Value = "  -2.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.21"
$ws.Range("E27").Value = "  +0.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.57"
$ws.Range("E28").Value = "  -1.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.018"
$ws.Range("E29").Value = "  -1.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.476"
$ws.Range("E30").Value = "  +3.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.615"
$ws.Range("E31").Value = "  +1.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.529"
$ws.Range("E32").Value = "  -1.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.294"
$ws.Range("E33").Value = "  +4.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05329"
$ws.Range("E34").Value = "  -3.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.234"
$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7528"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9998"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.696"
$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01954"
$ws.Range("E39").Value = "  -1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.750"
$ws.Range("E40").Value = "  -1.50%  "

$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.105.61"
$ws.Range("E42").Value = "  +1.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.080"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.19"
$ws.Range("E44").Value = "  -2.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8599"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.93"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.705"
$ws.Range("E48").Value = "  +0.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.842"
$ws.Range("E49").Value = "  -2.92%  "

$ws.Range("E50").Value = "  +1.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.014.74"
$ws.Range("E51").Value = "  -2.83%  "
